$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: fill in missing E18 (time out) ---
$ws.Range("E18").Value = 0.20833333333333334
$ws.Range("E18").NumberFormat = "h:mm"

# --- Row 19: complete the previously-empty row (Fri 13/07/2018) ---
$ws.Range("A19").Value = "7/13/2018"
$ws.Range("B19").Value = "N/A"
$ws.Range("B19").HorizontalAlignment = -4152
$ws.Range("C19").Value = "N/A"
$ws.Range("C19").HorizontalAlignment = -4152
$ws.Range("D19").Value = 0.375
$ws.Range("D19").NumberFormat = "h:mm"
$ws.Range("E19").Value = 0.20833333333333334
$ws.Range("E19").NumberFormat = "h:mm"

# --- Row 20: new entry, date stored as text "16/07/2018" ---
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "16/07/2018"
$ws.Range("B20").Value = "N/A"
$ws.Range("B20").HorizontalAlignment = -4152
$ws.Range("C20").Value = "N/A"
$ws.Range("C20").HorizontalAlignment = -4152
$ws.Range("D20").Value = 0.3125
$ws.Range("D20").NumberFormat = "h:mm"

# --- View: scroll back to top and move selection to E20 ---
$ws.Range("A1").Select()
$ws.Range("E20").Select()
